$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Columns E (Starttid) and F (Sluttid) on rows 3-7 were plain numbers before;
# they now need to match the text-formatted style already used on row 2
# (numFmtId 49 "@" plus matching borders), so copy the formatting down first.
$ws.Range("E2").Copy()
$ws.Range("E3:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy()
$ws.Range("F3:F7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update Starttid (E) / Sluttid (F) values for rows 2-7.
$ws.Range("E2").Value = "16:51"
$ws.Range("F2").Value = "18:01"

$ws.Range("E3").Value = "16:52"
$ws.Range("F3").Value = "18:02"

$ws.Range("E4").Value = "16:53"
$ws.Range("F4").Value = "18:03"

$ws.Range("E5").Value = "16:54"
$ws.Range("F5").Value = "18:04"

$ws.Range("E6").Value = "16:55"
$ws.Range("F6").Value = "18:05"

$ws.Range("E7").Value = "16:56"
$ws.Range("F7").Value = "18:01"

# Clear Q2, Q3, V2, V3 ("ti" / "sø" weekday flags no longer apply to those rows)
$ws.Range("Q2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("V3").ClearContents()

# Update the active selection to match the final state of the saved file.
$ws.Range("F7").Select()
